$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Home"
$ws.Range("C3").Value = 'SELECT  company_name,reg_id, company_address FROM public."Company"'
$ws.Range("D3").Value = "Company_1.csv"

$ws.Range("E3").Select()
